$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item("income").Name = "all_income"
$wb.Worksheets.Item("test").Name = "total_income"

# Update the "Sheet Name" label to "Opirations" on the total_income sheet (A1)
$ws4 = $wb.Worksheets.Item("total_income")
$ws4.Range("A1").Value = "Opirations"

# Update B2 value from 159217.9 to 0
$ws4.Range("B2").Value = 0
